$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as TEXT even when it looks numeric
# (mirrors a user formatting the cell as Text before typing the value, which
# is how Excel preserves strings such as "004634" or "0.50" verbatim).
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ===========================================================================
# 1) Build the new "2022-Q1" sheet by duplicating "2021-Q4" (so it inherits
#    the exact same column layout / cell styles), then overwrite its data.
# ===========================================================================
$q4 = $wb.Worksheets.Item(2)
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Drop the two trailing data rows copied over from "2021-Q4" (rows 5 and 6);
# the "2022-Q1" table only has 3 data rows.
$q1.Rows.Item(6).Delete()
$q1.Rows.Item(5).Delete()

# Row 2
Set-TextValue $q1 2 2 "004634"
$q1.Cells.Item(2,3).Value = "新疆前海联合泳涛灵活配置混合A"
Set-TextValue $q1 2 4 "1.33"
Set-TextValue $q1 2 5 "89.65"
Set-TextValue $q1 2 6 "4.85"
Set-TextValue $q1 2 7 "0.0645"
$q1.Cells.Item(2,8).Value = 6

# Row 3
Set-TextValue $q1 3 2 "006235"
$q1.Cells.Item(3,3).Value = "东方城镇消费主题混合"
Set-TextValue $q1 3 4 "0.50"
Set-TextValue $q1 3 5 "90.32"
Set-TextValue $q1 3 6 "4.81"
Set-TextValue $q1 3 7 "0.0240"
$q1.Cells.Item(3,8).Value = 5

# Row 4
Set-TextValue $q1 4 2 "007041"
$q1.Cells.Item(4,3).Value = "新疆前海联合泳涛灵活配置混合C"
Set-TextValue $q1 4 4 "0.00"
Set-TextValue $q1 4 5 "89.65"
Set-TextValue $q1 4 6 "4.85"
$q1.Cells.Item(4,7).Value = 0
$q1.Cells.Item(4,8).Value = 6

# ===========================================================================
# 2) Build the new "总计" (summary) sheet. Duplicate the existing "总计"
#    sheet (now pushed one slot to the right, after "2022-Q1") so the new
#    sheet inherits its layout/styles, append it at the end of the workbook,
#    delete the now-superseded original, and rename the duplicate.
# ===========================================================================
$oldSummary = $wb.Worksheets.Item(4)
$oldSummary.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item(4).Delete() | Out-Null
$summary = $wb.Worksheets.Item($wb.Worksheets.Count)
$summary.Name = "总计"

# Shift the existing two data rows down by one (bounded A:D range copy so we
# don't touch the rest of the (virtually infinite) row) to make room for the
# new "2022-Q1" row, carrying the existing cell styles along with them.
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# Row 2 becomes the new "2022-Q1" entry.
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q1"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.09

# Renumber the index column (A) for the rows that shifted down.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2

foreach ($s in $wb.Worksheets) {
    Write-Host $s.Index $s.Name
}
